$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")
$runtimes = $wb.Worksheets.Item("Runtimes")

# --- Data sheet: append rows 694-709 (A:O) ---
$data.Range("A694").Value = [double]"100"
$data.Range("B694").Value = [double]"5"
$data.Range("C694").Value = [double]"0.5"
$data.Range("D694").Value = 'ba-no-cycle'
$data.Range("E694").Value = '2021-06-17 11:37:32.720444'
$data.Range("F694").Value = [double]"4"
$data.Range("G694").Value = [double]"2"
$data.Range("H694").Value = [double]"2"
$data.Range("I694").Value = [double]"2"
$data.Range("J694").Value = [double]"15"
$data.Range("K694").Value = [double]"15"
$data.Range("L694").Value = [double]"15"
$data.Range("M694").Value = '-'
$data.Range("N694").Value = '-'
$data.Range("O694").Value = '-'

$data.Range("A695").Value = [double]"100"
$data.Range("B695").Value = [double]"5"
$data.Range("C695").Value = [double]"0.5"
$data.Range("D695").Value = 'ba-cycle'
$data.Range("E695").Value = '2021-06-17 11:37:34.256386'
$data.Range("F695").Value = '-'
$data.Range("G695").Value = '-'
$data.Range("H695").Value = '-'
$data.Range("I695").Value = '-'
$data.Range("J695").Value = '-'
$data.Range("K695").Value = '-'
$data.Range("L695").Value = '-'
$data.Range("M695").Value = [double]"15"
$data.Range("N695").Value = [double]"15"
$data.Range("O695").Value = [double]"15"

$data.Range("A696").Value = [double]"100"
$data.Range("B696").Value = [double]"5"
$data.Range("C696").Value = [double]"0.5"
$data.Range("D696").Value = 'er-no-cycle'
$data.Range("E696").Value = '2021-06-17 11:37:35.501490'
$data.Range("F696").Value = [double]"8"
$data.Range("G696").Value = [double]"9"
$data.Range("H696").Value = [double]"9"
$data.Range("I696").Value = [double]"9"
$data.Range("J696").Value = [double]"10"
$data.Range("K696").Value = [double]"10"
$data.Range("L696").Value = [double]"10"
$data.Range("M696").Value = '-'
$data.Range("N696").Value = '-'
$data.Range("O696").Value = '-'

$data.Range("A697").Value = [double]"100"
$data.Range("B697").Value = [double]"5"
$data.Range("C697").Value = [double]"0.5"
$data.Range("D697").Value = 'er-cycle'
$data.Range("E697").Value = '2021-06-17 11:37:36.580797'
$data.Range("F697").Value = '-'
$data.Range("G697").Value = '-'
$data.Range("H697").Value = '-'
$data.Range("I697").Value = '-'
$data.Range("J697").Value = '-'
$data.Range("K697").Value = '-'
$data.Range("L697").Value = '-'
$data.Range("M697").Value = [double]"10"
$data.Range("N697").Value = [double]"10"
$data.Range("O697").Value = [double]"10"

$data.Range("A698").Value = [double]"100"
$data.Range("B698").Value = [double]"5"
$data.Range("C698").Value = [double]"0.5"
$data.Range("D698").Value = 'ws-no-cycle'
$data.Range("E698").Value = '2021-06-17 11:37:37.738037'
$data.Range("F698").Value = [double]"1"
$data.Range("G698").Value = [double]"-4"
$data.Range("H698").Value = [double]"1"
$data.Range("I698").Value = [double]"1"
$data.Range("J698").Value = [double]"1"
$data.Range("K698").Value = [double]"1"
$data.Range("L698").Value = [double]"1"
$data.Range("M698").Value = '-'
$data.Range("N698").Value = '-'
$data.Range("O698").Value = '-'

$data.Range("A699").Value = [double]"100"
$data.Range("B699").Value = [double]"5"
$data.Range("C699").Value = [double]"0.5"
$data.Range("D699").Value = 'ws-cycle'
$data.Range("E699").Value = '2021-06-17 11:37:39.095243'
$data.Range("F699").Value = '-'
$data.Range("G699").Value = '-'
$data.Range("H699").Value = '-'
$data.Range("I699").Value = '-'
$data.Range("J699").Value = '-'
$data.Range("K699").Value = '-'
$data.Range("L699").Value = '-'
$data.Range("M699").Value = [double]"1"
$data.Range("N699").Value = [double]"1"
$data.Range("O699").Value = [double]"1"

$data.Range("A700").Value = [double]"100"
$data.Range("B700").Value = [double]"5"
$data.Range("C700").Value = [double]"0.5"
$data.Range("D700").Value = 'cluster no cycle'
$data.Range("E700").Value = '2021-06-17 11:37:40.352415'
$data.Range("F700").Value = [double]"14"
$data.Range("G700").Value = [double]"14"
$data.Range("H700").Value = [double]"11"
$data.Range("I700").Value = [double]"14"
$data.Range("J700").Value = [double]"14"
$data.Range("K700").Value = [double]"14"
$data.Range("L700").Value = [double]"14"
$data.Range("M700").Value = '-'
$data.Range("N700").Value = '-'
$data.Range("O700").Value = '-'

$data.Range("A701").Value = [double]"100"
$data.Range("B701").Value = [double]"5"
$data.Range("C701").Value = [double]"0.5"
$data.Range("D701").Value = 'cluster cycle'
$data.Range("E701").Value = '2021-06-17 11:37:41.731323'
$data.Range("F701").Value = '-'
$data.Range("G701").Value = '-'
$data.Range("H701").Value = '-'
$data.Range("I701").Value = '-'
$data.Range("J701").Value = '-'
$data.Range("K701").Value = '-'
$data.Range("L701").Value = '-'
$data.Range("M701").Value = [double]"13"
$data.Range("N701").Value = [double]"13"
$data.Range("O701").Value = [double]"13"

$data.Range("A702").Value = [double]"100"
$data.Range("B702").Value = [double]"5"
$data.Range("C702").Value = [double]"0.5"
$data.Range("D702").Value = 'ba-no-cycle'
$data.Range("E702").Value = '2021-06-17 11:38:15.167025'
$data.Range("F702").Value = [double]"16"
$data.Range("G702").Value = [double]"14"
$data.Range("H702").Value = [double]"14"
$data.Range("I702").Value = [double]"14"
$data.Range("J702").Value = [double]"14"
$data.Range("K702").Value = [double]"14"
$data.Range("L702").Value = [double]"14"
$data.Range("M702").Value = '-'
$data.Range("N702").Value = '-'
$data.Range("O702").Value = '-'

$data.Range("A703").Value = [double]"100"
$data.Range("B703").Value = [double]"5"
$data.Range("C703").Value = [double]"0.5"
$data.Range("D703").Value = 'ba-cycle'
$data.Range("E703").Value = '2021-06-17 11:38:16.363163'
$data.Range("F703").Value = '-'
$data.Range("G703").Value = '-'
$data.Range("H703").Value = '-'
$data.Range("I703").Value = '-'
$data.Range("J703").Value = '-'
$data.Range("K703").Value = '-'
$data.Range("L703").Value = '-'
$data.Range("M703").Value = [double]"14"
$data.Range("N703").Value = [double]"14"
$data.Range("O703").Value = [double]"14"

$data.Range("A704").Value = [double]"100"
$data.Range("B704").Value = [double]"5"
$data.Range("C704").Value = [double]"0.5"
$data.Range("D704").Value = 'er-no-cycle'
$data.Range("E704").Value = '2021-06-17 11:38:17.438332'
$data.Range("F704").Value = [double]"3"
$data.Range("G704").Value = [double]"-11"
$data.Range("H704").Value = [double]"0"
$data.Range("I704").Value = [double]"1"
$data.Range("J704").Value = [double]"1"
$data.Range("K704").Value = [double]"1"
$data.Range("L704").Value = [double]"-2"
$data.Range("M704").Value = '-'
$data.Range("N704").Value = '-'
$data.Range("O704").Value = '-'

$data.Range("A705").Value = [double]"100"
$data.Range("B705").Value = [double]"5"
$data.Range("C705").Value = [double]"0.5"
$data.Range("D705").Value = 'er-cycle'
$data.Range("E705").Value = '2021-06-17 11:38:18.664559'
$data.Range("F705").Value = '-'
$data.Range("G705").Value = '-'
$data.Range("H705").Value = '-'
$data.Range("I705").Value = '-'
$data.Range("J705").Value = '-'
$data.Range("K705").Value = '-'
$data.Range("L705").Value = '-'
$data.Range("M705").Value = [double]"1"
$data.Range("N705").Value = [double]"1"
$data.Range("O705").Value = [double]"-2"

$data.Range("A706").Value = [double]"100"
$data.Range("B706").Value = [double]"5"
$data.Range("C706").Value = [double]"0.5"
$data.Range("D706").Value = 'ws-no-cycle'
$data.Range("E706").Value = '2021-06-17 11:38:19.916148'
$data.Range("F706").Value = [double]"4"
$data.Range("G706").Value = [double]"2"
$data.Range("H706").Value = [double]"2"
$data.Range("I706").Value = [double]"2"
$data.Range("J706").Value = [double]"3"
$data.Range("K706").Value = [double]"0"
$data.Range("L706").Value = [double]"3"
$data.Range("M706").Value = '-'
$data.Range("N706").Value = '-'
$data.Range("O706").Value = '-'

$data.Range("A707").Value = [double]"100"
$data.Range("B707").Value = [double]"5"
$data.Range("C707").Value = [double]"0.5"
$data.Range("D707").Value = 'ws-cycle'
$data.Range("E707").Value = '2021-06-17 11:38:21.013486'
$data.Range("F707").Value = '-'
$data.Range("G707").Value = '-'
$data.Range("H707").Value = '-'
$data.Range("I707").Value = '-'
$data.Range("J707").Value = '-'
$data.Range("K707").Value = '-'
$data.Range("L707").Value = '-'
$data.Range("M707").Value = [double]"3"
$data.Range("N707").Value = [double]"0"
$data.Range("O707").Value = [double]"3"

$data.Range("A708").Value = [double]"100"
$data.Range("B708").Value = [double]"5"
$data.Range("C708").Value = [double]"0.5"
$data.Range("D708").Value = 'cluster no cycle'
$data.Range("E708").Value = '2021-06-17 11:38:22.083597'
$data.Range("F708").Value = [double]"14"
$data.Range("G708").Value = [double]"12"
$data.Range("H708").Value = [double]"13"
$data.Range("I708").Value = [double]"14"
$data.Range("J708").Value = [double]"14"
$data.Range("K708").Value = [double]"14"
$data.Range("L708").Value = [double]"13"
$data.Range("M708").Value = '-'
$data.Range("N708").Value = '-'
$data.Range("O708").Value = '-'

$data.Range("A709").Value = [double]"100"
$data.Range("B709").Value = [double]"5"
$data.Range("C709").Value = [double]"0.5"
$data.Range("D709").Value = 'cluster cycle'
$data.Range("E709").Value = '2021-06-17 11:38:23.418709'
$data.Range("F709").Value = '-'
$data.Range("G709").Value = '-'
$data.Range("H709").Value = '-'
$data.Range("I709").Value = '-'
$data.Range("J709").Value = '-'
$data.Range("K709").Value = '-'
$data.Range("L709").Value = '-'
$data.Range("M709").Value = [double]"4"
$data.Range("N709").Value = [double]"3"
$data.Range("O709").Value = [double]"3"

# --- Runtimes sheet: append rows 695-710 (A:L) ---
$runtimes.Range("A695").Value = [double]"100"
$runtimes.Range("B695").Value = [double]"5"
$runtimes.Range("C695").Value = [double]"0.5"
$runtimes.Range("D695").Value = 'ba-no-cycle'
$runtimes.Range("E695").Value = '2021-06-17 11:37:32.720444'
$runtimes.Range("F695").Value = [double]"0.04976590000000058"
$runtimes.Range("G695").Value = [double]"5.920000000081416e-05"
$runtimes.Range("H695").Value = [double]"0.04490209999999983"
$runtimes.Range("I695").Value = [double]"0.1264978000000001"
$runtimes.Range("J695").Value = [double]"0.03520579999999995"
$runtimes.Range("K695").Value = [double]"0.001519000000000048"
$runtimes.Range("L695").Value = [double]"0.004674200000000184"

$runtimes.Range("A696").Value = [double]"100"
$runtimes.Range("B696").Value = [double]"5"
$runtimes.Range("C696").Value = [double]"0.5"
$runtimes.Range("D696").Value = 'ba-cycle'
$runtimes.Range("E696").Value = '2021-06-17 11:37:34.256386'
$runtimes.Range("F696").Value = '-'
$runtimes.Range("G696").Value = '-'
$runtimes.Range("H696").Value = '-'
$runtimes.Range("I696").Value = '-'
$runtimes.Range("J696").Value = [double]"0.03518600000000038"
$runtimes.Range("K696").Value = [double]"0.0020876000000003"
$runtimes.Range("L696").Value = [double]"0.005920099999999984"

$runtimes.Range("A697").Value = [double]"100"
$runtimes.Range("B697").Value = [double]"5"
$runtimes.Range("C697").Value = [double]"0.5"
$runtimes.Range("D697").Value = 'er-no-cycle'
$runtimes.Range("E697").Value = '2021-06-17 11:37:35.501490'
$runtimes.Range("F697").Value = [double]"0.01122630000000058"
$runtimes.Range("G697").Value = [double]"4.139999999885902e-05"
$runtimes.Range("H697").Value = [double]"0.004135199999998562"
$runtimes.Range("I697").Value = [double]"0.02822910000000078"
$runtimes.Range("J697").Value = [double]"0.02811479999999911"
$runtimes.Range("K697").Value = [double]"0.0009005999999995851"
$runtimes.Range("L697").Value = [double]"0.002771700000000266"

$runtimes.Range("A698").Value = [double]"100"
$runtimes.Range("B698").Value = [double]"5"
$runtimes.Range("C698").Value = [double]"0.5"
$runtimes.Range("D698").Value = 'er-cycle'
$runtimes.Range("E698").Value = '2021-06-17 11:37:36.580797'
$runtimes.Range("F698").Value = '-'
$runtimes.Range("G698").Value = '-'
$runtimes.Range("H698").Value = '-'
$runtimes.Range("I698").Value = '-'
$runtimes.Range("J698").Value = [double]"0.02917229999999904"
$runtimes.Range("K698").Value = [double]"0.001201599999999914"
$runtimes.Range("L698").Value = [double]"0.003455799999999343"

$runtimes.Range("A699").Value = [double]"100"
$runtimes.Range("B699").Value = [double]"5"
$runtimes.Range("C699").Value = [double]"0.5"
$runtimes.Range("D699").Value = 'ws-no-cycle'
$runtimes.Range("E699").Value = '2021-06-17 11:37:37.738037'
$runtimes.Range("F699").Value = [double]"0.001196900000000056"
$runtimes.Range("G699").Value = [double]"7.62999999999181e-05"
$runtimes.Range("H699").Value = [double]"0.0008893000000007589"
$runtimes.Range("I699").Value = [double]"0.02579710000000013"
$runtimes.Range("J699").Value = [double]"0.02974460000000079"
$runtimes.Range("K699").Value = [double]"0.002602300000001279"
$runtimes.Range("L699").Value = [double]"0.004880999999999247"

$runtimes.Range("A700").Value = [double]"100"
$runtimes.Range("B700").Value = [double]"5"
$runtimes.Range("C700").Value = [double]"0.5"
$runtimes.Range("D700").Value = 'ws-cycle'
$runtimes.Range("E700").Value = '2021-06-17 11:37:39.095243'
$runtimes.Range("F700").Value = '-'
$runtimes.Range("G700").Value = '-'
$runtimes.Range("H700").Value = '-'
$runtimes.Range("I700").Value = '-'
$runtimes.Range("J700").Value = [double]"0.02962920000000047"
$runtimes.Range("K700").Value = [double]"0.002230700000000141"
$runtimes.Range("L700").Value = [double]"0.00433679999999903"

$runtimes.Range("A701").Value = [double]"100"
$runtimes.Range("B701").Value = [double]"5"
$runtimes.Range("C701").Value = [double]"0.5"
$runtimes.Range("D701").Value = 'cluster no cycle'
$runtimes.Range("E701").Value = '2021-06-17 11:37:40.352415'
$runtimes.Range("F701").Value = [double]"0.001701599999998749"
$runtimes.Range("G701").Value = [double]"6.899999999987472e-05"
$runtimes.Range("H701").Value = [double]"0.01722519999999861"
$runtimes.Range("I701").Value = [double]"0.03470490000000126"
$runtimes.Range("J701").Value = [double]"0.03040220000000105"
$runtimes.Range("K701").Value = [double]"0.001623300000000327"
$runtimes.Range("L701").Value = [double]"0.005242399999998426"

$runtimes.Range("A702").Value = [double]"100"
$runtimes.Range("B702").Value = [double]"5"
$runtimes.Range("C702").Value = [double]"0.5"
$runtimes.Range("D702").Value = 'cluster cycle'
$runtimes.Range("E702").Value = '2021-06-17 11:37:41.731323'
$runtimes.Range("F702").Value = '-'
$runtimes.Range("G702").Value = '-'
$runtimes.Range("H702").Value = '-'
$runtimes.Range("I702").Value = '-'
$runtimes.Range("J702").Value = [double]"0.03181549999999866"
$runtimes.Range("K702").Value = [double]"0.001281900000000391"
$runtimes.Range("L702").Value = [double]"0.004743299999999451"

$runtimes.Range("A703").Value = [double]"100"
$runtimes.Range("B703").Value = [double]"5"
$runtimes.Range("C703").Value = [double]"0.5"
$runtimes.Range("D703").Value = 'ba-no-cycle'
$runtimes.Range("E703").Value = '2021-06-17 11:38:15.167025'
$runtimes.Range("F703").Value = [double]"0.001415399999999956"
$runtimes.Range("G703").Value = [double]"4.319999999990998e-05"
$runtimes.Range("H703").Value = [double]"0.0006435000000002411"
$runtimes.Range("I703").Value = [double]"0.03166400000000014"
$runtimes.Range("J703").Value = [double]"0.03287930000000028"
$runtimes.Range("K703").Value = [double]"0.00128869999999992"
$runtimes.Range("L703").Value = [double]"0.003823500000000202"

$runtimes.Range("A704").Value = [double]"100"
$runtimes.Range("B704").Value = [double]"5"
$runtimes.Range("C704").Value = [double]"0.5"
$runtimes.Range("D704").Value = 'ba-cycle'
$runtimes.Range("E704").Value = '2021-06-17 11:38:16.363163'
$runtimes.Range("F704").Value = '-'
$runtimes.Range("G704").Value = '-'
$runtimes.Range("H704").Value = '-'
$runtimes.Range("I704").Value = '-'
$runtimes.Range("J704").Value = [double]"0.0306002000000003"
$runtimes.Range("K704").Value = [double]"0.001087700000000247"
$runtimes.Range("L704").Value = [double]"0.003167200000000037"

$runtimes.Range("A705").Value = [double]"100"
$runtimes.Range("B705").Value = [double]"5"
$runtimes.Range("C705").Value = [double]"0.5"
$runtimes.Range("D705").Value = 'er-no-cycle'
$runtimes.Range("E705").Value = '2021-06-17 11:38:17.438332'
$runtimes.Range("F705").Value = [double]"0.001273200000000863"
$runtimes.Range("G705").Value = [double]"4.739999999969768e-05"
$runtimes.Range("H705").Value = [double]"0.04664279999999987"
$runtimes.Range("I705").Value = [double]"0.03127879999999994"
$runtimes.Range("J705").Value = [double]"0.05611419999999967"
$runtimes.Range("K705").Value = [double]"0.002721799999999774"
$runtimes.Range("L705").Value = [double]"0.00522530000000021"

$runtimes.Range("A706").Value = [double]"100"
$runtimes.Range("B706").Value = [double]"5"
$runtimes.Range("C706").Value = [double]"0.5"
$runtimes.Range("D706").Value = 'er-cycle'
$runtimes.Range("E706").Value = '2021-06-17 11:38:18.664559'
$runtimes.Range("F706").Value = '-'
$runtimes.Range("G706").Value = '-'
$runtimes.Range("H706").Value = '-'
$runtimes.Range("I706").Value = '-'
$runtimes.Range("J706").Value = [double]"0.05648840000000011"
$runtimes.Range("K706").Value = [double]"0.002850500000000089"
$runtimes.Range("L706").Value = [double]"0.006328499999999515"

$runtimes.Range("A707").Value = [double]"100"
$runtimes.Range("B707").Value = [double]"5"
$runtimes.Range("C707").Value = [double]"0.5"
$runtimes.Range("D707").Value = 'ws-no-cycle'
$runtimes.Range("E707").Value = '2021-06-17 11:38:19.916148'
$runtimes.Range("F707").Value = [double]"0.0006739000000006712"
$runtimes.Range("G707").Value = [double]"2.939999999895804e-05"
$runtimes.Range("H707").Value = [double]"0.0002706000000003428"
$runtimes.Range("I707").Value = [double]"0.02959759999999889"
$runtimes.Range("J707").Value = [double]"0.03197600000000023"
$runtimes.Range("K707").Value = [double]"0.001204300000001268"
$runtimes.Range("L707").Value = [double]"0.003281199999999984"

$runtimes.Range("A708").Value = [double]"100"
$runtimes.Range("B708").Value = [double]"5"
$runtimes.Range("C708").Value = [double]"0.5"
$runtimes.Range("D708").Value = 'ws-cycle'
$runtimes.Range("E708").Value = '2021-06-17 11:38:21.013486'
$runtimes.Range("F708").Value = '-'
$runtimes.Range("G708").Value = '-'
$runtimes.Range("H708").Value = '-'
$runtimes.Range("I708").Value = '-'
$runtimes.Range("J708").Value = [double]"0.03321019999999919"
$runtimes.Range("K708").Value = [double]"0.001198799999999167"
$runtimes.Range("L708").Value = [double]"0.002681400000000167"

$runtimes.Range("A709").Value = [double]"100"
$runtimes.Range("B709").Value = [double]"5"
$runtimes.Range("C709").Value = [double]"0.5"
$runtimes.Range("D709").Value = 'cluster no cycle'
$runtimes.Range("E709").Value = '2021-06-17 11:38:22.083597'
$runtimes.Range("F709").Value = [double]"0.00184530000000116"
$runtimes.Range("G709").Value = [double]"5.100000000091143e-05"
$runtimes.Range("H709").Value = [double]"0.01588459999999969"
$runtimes.Range("I709").Value = [double]"0.03523429999999905"
$runtimes.Range("J709").Value = [double]"0.03513260000000074"
$runtimes.Range("K709").Value = [double]"0.002214999999999634"
$runtimes.Range("L709").Value = [double]"0.008323600000000653"

$runtimes.Range("A710").Value = [double]"100"
$runtimes.Range("B710").Value = [double]"5"
$runtimes.Range("C710").Value = [double]"0.5"
$runtimes.Range("D710").Value = 'cluster cycle'
$runtimes.Range("E710").Value = '2021-06-17 11:38:23.418709'
$runtimes.Range("F710").Value = '-'
$runtimes.Range("G710").Value = '-'
$runtimes.Range("H710").Value = '-'
$runtimes.Range("I710").Value = '-'
$runtimes.Range("J710").Value = [double]"0.0340246999999998"
$runtimes.Range("K710").Value = [double]"0.002088699999999832"
$runtimes.Range("L710").Value = [double]"0.007329300000000316"

# Clear the auto-applied column style so newly-created cells have no explicit `s` attribute (matches source, which never assigns per-cell style to data rows).
$data.Range("A694:O709").Style = "Normal"
$runtimes.Range("A695:L710").Style = "Normal"
